# Applies the "Final SRS Version 2.0" edits: selectively bolds key
# words/phrases across a handful of paragraphs (no wording changes,
# only w:b formatting added to specific sub-strings), and relocates
# the _GoBack bookmark.

$d = $word.ActiveDocument

function Set-BoldSubstring {
    param(
        [string]$Context,
        [string]$Target
    )

    $rng = $d.Content
    $ok = $rng.Find.Execute($Context)
    if (-not $ok) {
        Write-Host "CONTEXT NOT FOUND: $Context"
        return
    }

    $fullStart = $rng.Start
    $fullText = $rng.Text
    $idx = $fullText.IndexOf($Target)
    if ($idx -lt 0) {
        Write-Host "TARGET NOT IN CONTEXT: [$Target] within [$Context]"
        return
    }

    $boldStart = $fullStart + $idx
    $boldEnd = $boldStart + $Target.Length
    $boldRange = $d.Range($boldStart, $boldEnd)
    $boldRange.Font.Bold = $true
}

# --- 2.1.1 Customers paragraph -------------------------------------------
Set-BoldSubstring "Customers can browse different categories and find" "categories"
Set-BoldSubstring "customers can buy vouchers" "buy"
Set-BoldSubstring "Customers have access to all coupon" "all"
Set-BoldSubstring "customers is to register." "register"

# --- 2.1.2 Company Clients paragraph --------------------------------------
Set-BoldSubstring "lients are able to present their" "present"
Set-BoldSubstring "products to large customer" "large"
Set-BoldSubstring "customer audience which" "audience"
Set-BoldSubstring "viewed from everyone (even" "everyone"
Set-BoldSubstring "visitors which do not have registration" "do not"
Set-BoldSubstring "accessed and purchased only by registered" "only by"
Set-BoldSubstring "clients should apply and be" "apply"
Set-BoldSubstring "and be approved by the organization" "approved"
Set-BoldSubstring "the organization first. Otherwise" "first"
Set-BoldSubstring "offer will not be displayed" "not"

# --- System environment paragraph -----------------------------------------
Set-BoldSubstring "website through internet from" "internet"
Set-BoldSubstring "companies should register and fill" "register"
Set-BoldSubstring "fill in a request form and" "request"
Set-BoldSubstring "before they become partners which allows" "partners"

# Relocate the _GoBack bookmark from the end of the document to just after
# " on" in "...on the other hand have to register...".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$onRange = $d.Content
$onRange.Find.Execute(" on the other hand have to register") | Out-Null
$insertPos = $onRange.Start + " on".Length
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Set-BoldSubstring "but do not have to wait for approval" "not"
Set-BoldSubstring "but do not have to wait for approval" "wait"
Set-BoldSubstring "have to wait for approval. They" "approval"
Set-BoldSubstring "use the full potential of the website right away" "full potential"
Set-BoldSubstring "use the full potential of the website right away" "right away"

Write-Host "Done"
